# Read sample_text from configuration
#
# - frontend: add sample_text_en / sample_text_fr rows with EN/FR sample
#   message bodies used to pre-fill the text analysis.
# - id_labels: add a "Condition" column so each prefecture id/label pair can
#   be scoped to a department, plus three new rows (Antony, Boulogne-
#   Billancourt, Nanterre) for department 92.
# - text_analysis: flip "read_from_cache" (B9) to TRUE so the sample text is
#   now read from the configuration instead of being (re)generated.

$wb = $excel.ActiveWorkbook

$wsFrontend   = $wb.Worksheets.Item("frontend")
$wsBackend    = $wb.Worksheets.Item("backend")
$wsCaseFields = $wb.Worksheets.Item("case_fields")
$wsIdLabels   = $wb.Worksheets.Item("id_labels")
$wsTextAnal   = $wb.Worksheets.Item("text_analysis")

# ---------------------------------------------------------------------
# frontend: new rows 3 & 4 (sample_text_en / sample_text_fr)
# ---------------------------------------------------------------------

$frEn = "sample_text_en"
$frFr = "sample_text_fr"

$enText = "Hello," + [char]10 +
  "I am writing on behalf of one of our members, Mr C C, whose case is rather complex and whose residence" + [char]8209 + "permit renewal application number (num" + [char]233 + "ro de la demande de renouvellement de carte de s" + [char]233 + "jour) is 7500000000000000003." + [char]10 +
  "He currently holds subsidiary protection status and would like to travel to his home country, on a temporary and exceptional basis, to attend his father" + [char]8217 + "s funeral." + [char]10 +
  "He would also like to submit an asylum application." + [char]10 +
  "In addition, Mr C" + [char]8217 + "s processing" + [char]8209 + "extension certificate expired on 11 October 2024. He therefore wishes to obtain a new certificate so he can prove the legality of his stay while awaiting his residence permit." + [char]10 +
  "Without action in the next few days, he risks losing his job." + [char]10 +
  "Thank you in advance for your attention; please note the urgency" + [char]8212 + "his employment is at stake, so this matter is very important." + [char]10 +
  "Mr C would furthermore like to file an asylum request with France."

$frText = "Bonjour, " + [char]10 +
  "Je vous sollicite pour le compte de l'un de nos adh" + [char]233 + "rents, Monsieur C C, dont la situation est assez complexe, et dont le num" + [char]233 + "ro de la demande de renouvellement de carte de s" + [char]233 + "jour est le 7500000000000000003." + [char]10 +
  "Il b" + [char]233 + "n" + [char]233 + "ficie actuellement du statut de protection subsidiaire et aimerait se rendre dans son pays de fa" + [char]231 + "on temporaire et exceptionnelle pour assister aux obs" + [char]232 + "ques de son p" + [char]232 + "re." + [char]10 +
  "Il aimerait d" + [char]233 + "poser une demande de droit d'asile." + [char]10 +
  "Par ailleurs, l'attestation de prolongation d'instruction de Monsieur C est arriv" + [char]233 + "e " + [char]224 + " expiration depuis le 11 octobre 2024. Aussi, il souhaiterait obtenir une nouvelle attestation pour pouvoir justifier de la r" + [char]233 + "gularit" + [char]233 + " de son s" + [char]233 + "jour, dans l'attente de recevoir carte de s" + [char]233 + "jour. " + [char]10 +
  "Sans action dans les prochains jours, il risquera de perdre son travail." + [char]10 +
  "Je vous remercie par avance et vous prie de noter l'urgence. Il risque son emploi, c'est donc tr" + [char]232 + "s important." + [char]10 +
  "Monsieur C aimerait, par ailleurs, faire une demande d'Asile " + [char]224 + " la France."

# Row 3: sample_text_en / EN text -- reuse existing formatting (A2 -> A/C style, C2 -> B style)
$wsFrontend.Range("A2").Copy()
$wsFrontend.Range("A3").PasteSpecial(-4122)
$wsFrontend.Range("C2").Copy()
$wsFrontend.Range("B3").PasteSpecial(-4122)
$wsFrontend.Range("A2").Copy()
$wsFrontend.Range("C3").PasteSpecial(-4122)
$wsFrontend.Range("A3").Value = $frEn
$wsFrontend.Range("B3").Value = $enText
$wsFrontend.Rows.Item(3).RowHeight = 255

# Row 4: sample_text_fr / FR text
$wsFrontend.Range("A2").Copy()
$wsFrontend.Range("A4").PasteSpecial(-4122)
$wsFrontend.Range("C2").Copy()
$wsFrontend.Range("B4").PasteSpecial(-4122)
$wsFrontend.Range("A2").Copy()
$wsFrontend.Range("C4").PasteSpecial(-4122)
$wsFrontend.Range("A4").Value = $frFr
$wsFrontend.Range("B4").Value = $frText
$wsFrontend.Rows.Item(4).RowHeight = 270

$wsFrontend.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# id_labels: new "Condition" column (D) + new rows for dpt 92
# ---------------------------------------------------------------------

$cond78 = "{departement} == 78"
$cond92 = "{departement} == 92"

# Header D1, formatted like A1
$wsIdLabels.Range("A1").Copy()
$wsIdLabels.Range("D1").PasteSpecial(-4122)
$wsIdLabels.Range("D1").Value = "Condition"

# D2:D5, formatted like A2, all condition = dpt 78
$wsIdLabels.Range("A2").Copy()
$wsIdLabels.Range("D2:D5").PasteSpecial(-4122)
$wsIdLabels.Range("D2").Value = $cond78
$wsIdLabels.Range("D3").Value = $cond78
$wsIdLabels.Range("D4").Value = $cond78
$wsIdLabels.Range("D5").Value = $cond78

# New rows 6-8, formatted like row 5
$wsIdLabels.Range("A5:D5").Copy()
$wsIdLabels.Range("A6:D8").PasteSpecial(-4122)

$wsIdLabels.Range("A6").Value = "ANT"
$wsIdLabels.Range("B6").Value = "Antony"
$wsIdLabels.Range("C6").Value = "Antony"
$wsIdLabels.Range("D6").Value = $cond92

$wsIdLabels.Range("A7").Value = "BOULO"
$wsIdLabels.Range("B7").Value = "Boulogne-Billancourt"
$wsIdLabels.Range("C7").Value = "Boulogne-Billancourt"
$wsIdLabels.Range("D7").Value = $cond92

$wsIdLabels.Range("A8").Value = "NANTER"
$wsIdLabels.Range("B8").Value = "Nanterre"
$wsIdLabels.Range("C8").Value = "Nanterre"
$wsIdLabels.Range("D8").Value = $cond92

$wsIdLabels.Columns.Item(4).ColumnWidth = 18.3

$wsIdLabels.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# text_analysis: read_from_cache (B9) False -> True
# ---------------------------------------------------------------------

$wsTextAnal.Range("B9").Value = $true

# ---------------------------------------------------------------------
# Selections / active sheet bookkeeping (matches the recorded UI state)
# ---------------------------------------------------------------------

$wsTextAnal.Range("B9").Select()
$wsBackend.Range("C7").Select()
$wsCaseFields.Range("A3").Select()
$wsIdLabels.Range("E6").Select()
$wsFrontend.Activate()
$wsFrontend.Range("B4").Select()
